$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 36.71718733333334
$ws.Range("H2").Value = 110.151562
$ws.Range("I2").Value = 0.728743057485239
$ws.Range("J2").Value = 0.7328478088626956
$ws.Range("M2").Value = 0.4349816666666667
$ws.Range("N2").Value = 1.304945
$ws.Range("O2").Value = 0.0068619340787224
$ws.Range("P2").Value = 0.00687614195861717
$ws.Range("Q2").Value = 15.97130334156556
$ws.Range("R2").Value = 143.74173007409
$ws.Range("S2").Value = 0.005000586820790319
$ws.Range("T2").Value = 0.005039165567801437

$ws.Range("G3").Value = 36.71718733333334
$ws.Range("H3").Value = 110.151562
$ws.Range("I3").Value = 0.728743057485239
$ws.Range("J3").Value = 0.7328478088626956
$ws.Range("O3").Value = 0.005722841821244052
$ws.Range("P3").Value = 0.00573469116988544
$ws.Range("Q3").Value = 13.32004091766267
$ws.Range("R3").Value = 119.880368258964
$ws.Range("S3").Value = 0.004170481246317784
$ws.Range("T3").Value = 0.004202655858354793

$ws.Range("G4").Value = 36.71718733333334
$ws.Range("H4").Value = 110.151562
$ws.Range("I4").Value = 0.728743057485239
$ws.Range("J4").Value = 0.7328478088626956
$ws.Range("M4").Value = 29.338587
$ws.Range("N4").Value = 88.015761
$ws.Range("O4").Value = 0.4628228391775791
$ws.Range("P4").Value = 0.4637811304167767
$ws.Range("Q4").Value = 1077.230394974298
$ws.Range("R4").Value = 9695.073554768682
$ws.Range("S4").Value = 0.3372789308962681
$ws.Range("T4").Value = 0.3398809852177989

$ws.Range("G5").Value = 36.71718733333334
$ws.Range("H5").Value = 110.151562
$ws.Range("I5").Value = 0.728743057485239
$ws.Range("J5").Value = 0.7328478088626956
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.3929435
$ws.Range("N5").Value = 0.785887
$ws.Range("O5").Value = 0.006198772500747056
$ws.Range("P5").Value = 0.004141071520586516
$ws.Range("Q5").Value = 14.42778010091567
$ws.Range("R5").Value = 86.566680605494
$ws.Range("S5").Value = 0.004517312424849831
$ws.Range("T5").Value = 0.00303477519020554

$ws.Range("G6").Value = 36.71718733333334
$ws.Range("H6").Value = 110.151562
$ws.Range("I6").Value = 0.728743057485239
$ws.Range("J6").Value = 0.7328478088626956
$ws.Range("M6").Value = 32.861248
$ws.Range("N6").Value = 98.583744
$ws.Range("O6").Value = 0.5183936124217073
$ws.Range("P6").Value = 0.5194669649341341
$ws.Range("Q6").Value = 1206.572598823125
$ws.Range("R6").Value = 10859.15338940813
$ws.Range("S6").Value = 0.377775746097013
$ws.Range("T6").Value = 0.3806902270285349

$ws.Range("I7").Value = 0.02522574977045663
$ws.Range("J7").Value = 0.0253678374789488
$ws.Range("M7").Value = 0.4349816666666667
$ws.Range("N7").Value = 1.304945
$ws.Range("O7").Value = 0.0068619340787224
$ws.Range("P7").Value = 0.00687614195861717
$ws.Range("Q7").Value = 0.5528534336816666
$ws.Range("R7").Value = 4.975680903134999
$ws.Range("S7").Value = 0.0001730974320112202
$ws.Range("T7").Value = 0.0001744328516883811

$ws.Range("I8").Value = 0.02522574977045663
$ws.Range("J8").Value = 0.0253678374789488
$ws.Range("O8").Value = 0.005722841821244052
$ws.Range("P8").Value = 0.00573469116988544
$ws.Range("S8").Value = 0.0001443629757586068
$ws.Range("T8").Value = 0.0001454767135896166

$ws.Range("I9").Value = 0.02522574977045663
$ws.Range("J9").Value = 0.0253678374789488
$ws.Range("M9").Value = 29.338587
$ws.Range("N9").Value = 88.015761
$ws.Range("O9").Value = 0.4628228391775791
$ws.Range("P9").Value = 0.4637811304167767
$ws.Range("Q9").Value = 37.288786643847
$ws.Range("R9").Value = 335.599079794623
$ws.Range("S9").Value = 0.0116750531291459
$ws.Range("T9").Value = 0.01176512434221595

$ws.Range("I10").Value = 0.02522574977045663
$ws.Range("J10").Value = 0.0253678374789488
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.5
$ws.Range("M10").Value = 0.3929435
$ws.Range("N10").Value = 0.785887
$ws.Range("O10").Value = 0.006198772500747056
$ws.Range("P10").Value = 0.004141071520586516
$ws.Range("Q10").Value = 0.4994237225735
$ws.Range("R10").Value = 2.996542335441
$ws.Range("S10").Value = 0.0001563686839878329
$ws.Range("T10").Value = 0.0001050500293229421

$ws.Range("I11").Value = 0.02522574977045663
$ws.Range("J11").Value = 0.0253678374789488
$ws.Range("M11").Value = 32.861248
$ws.Range("N11").Value = 98.583744
$ws.Range("O11").Value = 0.5183936124217073
$ws.Range("P11").Value = 0.5194669649341341
$ws.Range("Q11").Value = 41.76602184428799
$ws.Range("R11").Value = 375.894196598592
$ws.Range("S11").Value = 0.01307686754955307
$ws.Range("T11").Value = 0.01317775354213191

$ws.Range("G12").Value = 7.275657333333332
$ws.Range("H12").Value = 21.826972
$ws.Range("I12").Value = 0.1444033477339586
$ws.Range("J12").Value = 0.1452167206154317
$ws.Range("M12").Value = 0.4349816666666667
$ws.Range("N12").Value = 1.304945
$ws.Range("O12").Value = 0.0068619340787224
$ws.Range("P12").Value = 0.00687614195861717
$ws.Range("Q12").Value = 3.164777552948888
$ws.Range("R12").Value = 28.48299797654
$ws.Range("S12").Value = 0.0009908862528972516
$ws.Range("T12").Value = 0.000998530785716557

$ws.Range("G13").Value = 7.275657333333332
$ws.Range("H13").Value = 21.826972
$ws.Range("I13").Value = 0.1444033477339586
$ws.Range("J13").Value = 0.1452167206154317
$ws.Range("O13").Value = 0.005722841821244052
$ws.Range("P13").Value = 0.00573469116988544
$ws.Range("Q13").Value = 2.639419313442666
$ws.Range("R13").Value = 23.754773820984
$ws.Range("S13").Value = 0.0008263975175395458
$ws.Range("T13").Value = 0.0008327730454330373

$ws.Range("G14").Value = 7.275657333333332
$ws.Range("H14").Value = 21.826972
$ws.Range("I14").Value = 0.1444033477339586
$ws.Range("J14").Value = 0.1452167206154317
$ws.Range("M14").Value = 29.338587
$ws.Range("N14").Value = 88.015761
$ws.Range("O14").Value = 0.4628228391775791
$ws.Range("P14").Value = 0.4637811304167767
$ws.Range("Q14").Value = 213.457505656188
$ws.Range("R14").Value = 1921.117550905692
$ws.Range("S14").Value = 0.06683316738497795
$ws.Range("T14").Value = 0.06734877484244217

$ws.Range("G15").Value = 7.275657333333332
$ws.Range("H15").Value = 21.826972
$ws.Range("I15").Value = 0.1444033477339586
$ws.Range("J15").Value = 0.1452167206154317
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.5
$ws.Range("M15").Value = 0.3929435
$ws.Range("N15").Value = 0.785887
$ws.Range("O15").Value = 0.006198772500747056
$ws.Range("P15").Value = 0.004141071520586516
$ws.Range("Q15").Value = 2.858922257360666
$ws.Range("R15").Value = 17.153533544164
$ws.Range("S15").Value = 0.0008951235009490773
$ws.Range("T15").Value = 0.0006013528260535332

$ws.Range("G16").Value = 7.275657333333332
$ws.Range("H16").Value = 21.826972
$ws.Range("I16").Value = 0.1444033477339586
$ws.Range("J16").Value = 0.1452167206154317
$ws.Range("M16").Value = 32.861248
$ws.Range("N16").Value = 98.583744
$ws.Range("O16").Value = 0.5183936124217073
$ws.Range("P16").Value = 0.5194669649341341
$ws.Range("Q16").Value = 239.0871799936853
$ws.Range("R16").Value = 2151.784619943168
$ws.Range("S16").Value = 0.07485777307759477
$ws.Range("T16").Value = 0.07543528911578642

$ws.Range("G17").Value = 0.8466215
$ws.Range("H17").Value = 1.693243
$ws.Range("I17").Value = 0.01680328982832053
$ws.Range("J17").Value = 0.0112652912032432
$ws.Range("M17").Value = 0.4349816666666667
$ws.Range("N17").Value = 1.304945
$ws.Range("O17").Value = 0.0068619340787224
$ws.Range("P17").Value = 0.00687614195861717
$ws.Range("Q17").Value = 0.3682648311058334
$ws.Range("R17").Value = 2.209588986635
$ws.Range("S17").Value = 0.0001153030671076021
$ws.Range("T17").Value = 0.00007746174151866144

$ws.Range("G18").Value = 0.8466215
$ws.Range("H18").Value = 1.693243
$ws.Range("I18").Value = 0.01680328982832053
$ws.Range("J18").Value = 0.0112652912032432
$ws.Range("O18").Value = 0.005722841821244052
$ws.Range("P18").Value = 0.00573469116988544
$ws.Range("Q18").Value = 0.307132268041
$ws.Range("R18").Value = 1.842793608246
$ws.Range("S18").Value = 0.00009616256976399752
$ws.Range("T18").Value = 0.00006460296598942688

$ws.Range("G19").Value = 0.8466215
$ws.Range("H19").Value = 1.693243
$ws.Range("I19").Value = 0.01680328982832053
$ws.Range("J19").Value = 0.0112652912032432
$ws.Range("M19").Value = 29.338587
$ws.Range("N19").Value = 88.015761
$ws.Range("O19").Value = 0.4628228391775791
$ws.Range("P19").Value = 0.4637811304167767
$ws.Range("Q19").Value = 24.8386785338205
$ws.Range("R19").Value = 149.032071202923
$ws.Range("S19").Value = 0.007776946305867043
$ws.Range("T19").Value = 0.0052246294887143

$ws.Range("G20").Value = 0.8466215
$ws.Range("H20").Value = 1.693243
$ws.Range("I20").Value = 0.01680328982832053
$ws.Range("J20").Value = 0.0112652912032432
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.5
$ws.Range("M20").Value = 0.3929435
$ws.Range("N20").Value = 0.785887
$ws.Range("O20").Value = 0.006198772500747056
$ws.Range("P20").Value = 0.004141071520586516
$ws.Range("Q20").Value = 0.33267441538525
$ws.Range("R20").Value = 1.330697661541
$ws.Range("S20").Value = 0.000104159770909876
$ws.Range("T20").Value = 0.0000466503765728642

$ws.Range("G21").Value = 0.8466215
$ws.Range("H21").Value = 1.693243
$ws.Range("I21").Value = 0.01680328982832053
$ws.Range("J21").Value = 0.0112652912032432
$ws.Range("M21").Value = 32.861248
$ws.Range("N21").Value = 98.583744
$ws.Range("O21").Value = 0.5183936124217073
$ws.Range("P21").Value = 0.5194669649341341
$ws.Range("Q21").Value = 27.821039073632
$ws.Range("R21").Value = 166.926234441792
$ws.Range("S21").Value = 0.008710718114672009
$ws.Range("T21").Value = 0.005851946630447942

$ws.Range("G22").Value = 4.273823333333334
$ws.Range("H22").Value = 12.82147
$ws.Range("I22").Value = 0.08482455518202518
$ws.Range("J22").Value = 0.08530234183968073
$ws.Range("M22").Value = 0.4349816666666667
$ws.Range("N22").Value = 1.304945
$ws.Range("O22").Value = 0.0068619340787224
$ws.Range("P22").Value = 0.00687614195861717
$ws.Range("Q22").Value = 1.859034796572222
$ws.Range("R22").Value = 16.73131316915
$ws.Range("S22").Value = 0.0005820605059160074
$ws.Range("T22").Value = 0.0005865510118921336

$ws.Range("G23").Value = 4.273823333333334
$ws.Range("H23").Value = 12.82147
$ws.Range("I23").Value = 0.08482455518202518
$ws.Range("J23").Value = 0.08530234183968073
$ws.Range("O23").Value = 0.005722841821244052
$ws.Range("P23").Value = 0.00573469116988544
$ws.Range("Q23").Value = 1.550431985926667
$ws.Range("R23").Value = 13.95388787334
$ws.Range("S23").Value = 0.0004854375118641176
$ws.Range("T23").Value = 0.0004891825865185664

$ws.Range("G24").Value = 4.273823333333334
$ws.Range("H24").Value = 12.82147
$ws.Range("I24").Value = 0.08482455518202518
$ws.Range("J24").Value = 0.08530234183968073
$ws.Range("M24").Value = 29.338587
$ws.Range("N24").Value = 88.015761
$ws.Range("O24").Value = 0.4628228391775791
$ws.Range("P24").Value = 0.4637811304167767
$ws.Range("Q24").Value = 125.38793768763
$ws.Range("R24").Value = 1128.49143918867
$ws.Range("S24").Value = 0.03925874146132013
$ws.Range("T24").Value = 0.03956161652560544

$ws.Range("G25").Value = 4.273823333333334
$ws.Range("H25").Value = 12.82147
$ws.Range("I25").Value = 0.08482455518202518
$ws.Range("J25").Value = 0.08530234183968073
$ws.Range("K25").Value = 1
$ws.Range("L25").Value = 0.5
$ws.Range("M25").Value = 0.3929435
$ws.Range("N25").Value = 0.785887
$ws.Range("O25").Value = 0.006198772500747056
$ws.Range("P25").Value = 0.004141071520586516
$ws.Range("Q25").Value = 1.679371098981667
$ws.Range("R25").Value = 10.07622659389
$ws.Range("S25").Value = 0.0005258081200504389
$ws.Range("T25").Value = 0.0003532430984316375

$ws.Range("G26").Value = 4.273823333333334
$ws.Range("H26").Value = 12.82147
$ws.Range("I26").Value = 0.08482455518202518
$ws.Range("J26").Value = 0.08530234183968073
$ws.Range("M26").Value = 32.861248
$ws.Range("N26").Value = 98.583744
$ws.Range("O26").Value = 0.5183936124217073
$ws.Range("P26").Value = 0.5194669649341341
$ws.Range("Q26").Value = 140.4431684648533
$ws.Range("R26").Value = 1263.98851618368
$ws.Range("S26").Value = 0.04397250758287449
$ws.Range("T26").Value = 0.04431174861723295
